# decomptes_DE.xlsx - "DE data now (finally) clean"
# The "Données" sheet's per-region date-range breakdown was re-derived:
# region "00" (UNIQUE) gains a missing row (C=2 / 40884-41784), and every
# row below it shifts down by one. A stray format-only fill (date number
# format) was also left behind in columns L:M for the first several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Données")
$limits = $wb.Worksheets.Item("Limites")

# --- Extend formatting of the data block down to the new last row (20),
#     mirroring row 19's per-cell styles (text-format col A is already
#     covered by the column style; D/E need the date number format).
$ws.Range("D19:E19").Copy()
$ws.Range("D20:E20").PasteSpecial(-4122)

# --- Rewrite the 19 data rows (2-20) with the corrected values ---
$ws.Cells.Item(2, 1).Value = "00"
$ws.Cells.Item(2, 2).Value = "UNIQUE"
$ws.Cells.Item(2, 3).Value = 81
$ws.Cells.Item(2, 4).Value = 29013
$ws.Cells.Item(2, 5).Value = 34496

$ws.Cells.Item(3, 1).Value = "00"
$ws.Cells.Item(3, 2).Value = "UNIQUE"
$ws.Cells.Item(3, 3).Value = 87
$ws.Cells.Item(3, 4).Value = 34497
$ws.Cells.Item(3, 5).Value = 38150

$ws.Cells.Item(4, 1).Value = "00"
$ws.Cells.Item(4, 2).Value = "UNIQUE"
$ws.Cells.Item(4, 3).Value = 74
$ws.Cells.Item(4, 4).Value = 43611
$ws.Cells.Item(4, 5).Value = 43861

$ws.Cells.Item(5, 1).Value = "00"
$ws.Cells.Item(5, 2).Value = "UNIQUE"
$ws.Cells.Item(5, 3).Value = 2
$ws.Cells.Item(5, 4).Value = 40884
$ws.Cells.Item(5, 5).Value = 41784
# E5 is a brand-new cell - give it the same date format as the rest of column E
$ws.Range("D2").Copy()
$ws.Range("E5").PasteSpecial(-4122)

$ws.Cells.Item(6, 1).Value = "00"
$ws.Cells.Item(6, 2).Value = "UNIQUE"
$ws.Cells.Item(6, 3).Value = 79
$ws.Cells.Item(6, 4).Value = 43862
$ws.Cells.Item(6, 5).Value = "NA"
# E6 used to carry the date format (it held a real date); now that it's
# the "NA" text marker it must drop back to the unstyled default
$ws.Range("B2").Copy()
$ws.Range("E6").PasteSpecial(-4122)

$ws.Cells.Item(7, 1).Value = "03"
$ws.Cells.Item(7, 2).Value = "EST"
$ws.Cells.Item(7, 3).Value = 10
$ws.Cells.Item(7, 4).Value = 38151
$ws.Cells.Item(7, 5).Value = 39970

$ws.Cells.Item(8, 1).Value = "03"
$ws.Cells.Item(8, 2).Value = "EST"
$ws.Cells.Item(8, 3).Value = 9
$ws.Cells.Item(8, 4).Value = 39971
$ws.Cells.Item(8, 5).Value = 43610

$ws.Cells.Item(9, 1).Value = "07"
$ws.Cells.Item(9, 2).Value = "ILE DE FRANCE"
$ws.Cells.Item(9, 3).Value = 14
$ws.Cells.Item(9, 4).Value = 38151
$ws.Cells.Item(9, 5).Value = 39970

$ws.Cells.Item(10, 1).Value = "07"
$ws.Cells.Item(10, 2).Value = "ILE DE FRANCE"
$ws.Cells.Item(10, 3).Value = 13
$ws.Cells.Item(10, 4).Value = 39971
$ws.Cells.Item(10, 5).Value = 41783

$ws.Cells.Item(11, 1).Value = "07"
$ws.Cells.Item(11, 2).Value = "ILE DE FRANCE"
$ws.Cells.Item(11, 3).Value = 15
$ws.Cells.Item(11, 4).Value = 41784
$ws.Cells.Item(11, 5).Value = 43610

$ws.Cells.Item(12, 1).Value = "06"
$ws.Cells.Item(12, 2).Value = "MASSIF CENTRAL CENTRE"
$ws.Cells.Item(12, 3).Value = 6
$ws.Cells.Item(12, 4).Value = 38151
$ws.Cells.Item(12, 5).Value = 39970

$ws.Cells.Item(13, 1).Value = "06"
$ws.Cells.Item(13, 2).Value = "MASSIF CENTRAL CENTRE"
$ws.Cells.Item(13, 3).Value = 5
$ws.Cells.Item(13, 4).Value = 39971
$ws.Cells.Item(13, 5).Value = 43610

$ws.Cells.Item(14, 1).Value = "01"
$ws.Cells.Item(14, 2).Value = "NORD OUEST"
$ws.Cells.Item(14, 3).Value = 12
$ws.Cells.Item(14, 4).Value = 38151
$ws.Cells.Item(14, 5).Value = 39970

$ws.Cells.Item(15, 1).Value = "01"
$ws.Cells.Item(15, 2).Value = "NORD OUEST"
$ws.Cells.Item(15, 3).Value = 10
$ws.Cells.Item(15, 4).Value = 39971
$ws.Cells.Item(15, 5).Value = 43610

$ws.Cells.Item(16, 1).Value = "02"
$ws.Cells.Item(16, 2).Value = "OUEST"
$ws.Cells.Item(16, 3).Value = 10
$ws.Cells.Item(16, 4).Value = 38151
$ws.Cells.Item(16, 5).Value = 39970

$ws.Cells.Item(17, 1).Value = "02"
$ws.Cells.Item(17, 2).Value = "OUEST"
$ws.Cells.Item(17, 3).Value = 9
$ws.Cells.Item(17, 4).Value = 39971
$ws.Cells.Item(17, 5).Value = 43610

$ws.Cells.Item(18, 1).Value = "08"
$ws.Cells.Item(18, 2).Value = "OUTRE MER"
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 4).Value = 38151
$ws.Cells.Item(18, 5).Value = 43610

$ws.Cells.Item(19, 1).Value = "05"
$ws.Cells.Item(19, 2).Value = "SUD EST"
$ws.Cells.Item(19, 3).Value = 13
$ws.Cells.Item(19, 4).Value = 38151
$ws.Cells.Item(19, 5).Value = 43610

$ws.Cells.Item(20, 1).Value = "04"
$ws.Cells.Item(20, 2).Value = "SUD OUEST"
$ws.Cells.Item(20, 3).Value = 10
$ws.Cells.Item(20, 4).Value = 38151
$ws.Cells.Item(20, 5).Value = 43610

# --- Stray date-format (no value) fill left in columns L:M for the first
#     rows of the block, as in the source edit ---
$ws.Range("D2").Copy()
$ws.Range("L2:M6").PasteSpecial(-4122)
$ws.Range("L7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selections as left by the editor ---
$ws.Range("E5").Select()
$limits.Range("A2").Select()
